# Auto-generated edit script applying the cell-level diffs to the cryptos worksheet.
# Updates columns B/C/D/E for rows 2-51 to reflect the refreshed crypto price feed:
#  - most rows keep the same coin but get new Price/Volume(1h) values
#  - rows 37-51 drop 'Frax' and shift the remaining coins up by one, appending 'Decentraland'
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = '@'
    $rng.Value = $value
    $rng.Style = 'Normal'
}

$ws.Range('D2').Value = '30.416.18'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = '1.927.43'
$ws.Range('E3').Value = '  +4.09%  '
Set-TextValue 'D4' '0.9992'
$ws.Range('E4').Value = '  -0.09%  '
Set-TextValue 'D5' '240.09'
$ws.Range('E5').Value = '  +2.87%  '
Set-TextValue 'D6' '0.9994'
$ws.Range('E6').Value = '  -0.08%  '
Set-TextValue 'D7' '0.4760'
$ws.Range('E7').Value = '  +0.41%  '
Set-TextValue 'D8' '0.2870'
$ws.Range('E8').Value = '  +4.55%  '
Set-TextValue 'D9' '0.06583'
$ws.Range('E9').Value = '  +4.01%  '
Set-TextValue 'D10' '19.12'
$ws.Range('E10').Value = '  +7.91%  '
Set-TextValue 'D11' '106.95'
$ws.Range('E11').Value = '  +26.33%  '
$ws.Range('D12').Value = '1.923.71'
$ws.Range('E12').Value = '  +3.80%  '
$ws.Range('E13').Value = '  +2.36%  '
Set-TextValue 'D14' '5.147'
$ws.Range('E14').Value = '  +3.96%  '
Set-TextValue 'D15' '0.6596'
Set-TextValue 'D16' '309.04'
$ws.Range('E16').Value = '  +25.88%  '
$ws.Range('D17').Value = '30.435.20'
$ws.Range('E17').Value = '  +0.14%  '
Set-TextValue 'D18' '13.01'
$ws.Range('E18').Value = '  +2.83%  '
Set-TextValue 'D19' '0.9991'
$ws.Range('E19').Value = '  -0.10%  '
Set-TextValue 'D20' '0.000007529'
$ws.Range('E20').Value = '  +2.72%  '
$ws.Range('D21').Value = '2.172.40'
$ws.Range('E21').Value = '  +3.26%  '
Set-TextValue 'D22' '5.323'
$ws.Range('E22').Value = '  +8.68%  '
Set-TextValue 'D23' '0.9987'
$ws.Range('E23').Value = '  -0.28%  '
Set-TextValue 'D24' '6.293'
$ws.Range('E24').Value = '  +6.66%  '
Set-TextValue 'D25' '167.39'
$ws.Range('E25').Value = '  +1.54%  '
Set-TextValue 'D26' '9.251'
$ws.Range('E26').Value = '  +2.03%  '
Set-TextValue 'D27' '20.27'
$ws.Range('E27').Value = '  +12.81%  '
$ws.Range('E28').Value = '  +8.78%  '
Set-TextValue 'D29' '0.1113'
$ws.Range('E29').Value = '  +8.42%  '
Set-TextValue 'D30' '1.359'
$ws.Range('E30').Value = '  +0.84%  '
Set-TextValue 'D31' '4.106'
$ws.Range('E31').Value = '  +1.68%  '
Set-TextValue 'D32' '3.934'
$ws.Range('E32').Value = '  +3.05%  '
Set-TextValue 'D33' '0.05019'
$ws.Range('E33').Value = '  +3.96%  '
Set-TextValue 'D34' '0.7439'
$ws.Range('E34').Value = '  +6.72%  '
Set-TextValue 'D35' '1.149'
$ws.Range('E35').Value = '  +1.90%  '
Set-TextValue 'D36' '2.760'
$ws.Range('E36').Value = '  +1.84%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D37' '0.01961'
$ws.Range('E37').Value = '  +3.25%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D38' '2.708'
$ws.Range('E38').Value = '  +0.88%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D39' '2.046'
$ws.Range('E39').Value = '  +2.90%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D40' '0.8784'
$ws.Range('E40').Value = '  +0.06%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D41' '106.99'
$ws.Range('E41').Value = '  +0.47%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D42' '70.93'
$ws.Range('E42').Value = '  +12.34%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D43' '5.806'
$ws.Range('E43').Value = '  +5.49%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D44' '0.9990'
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D45' '0.4173'
$ws.Range('E45').Value = '  +2.83%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D46' '7.276'
$ws.Range('E46').Value = '  +1.51%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D47' '9.280'
$ws.Range('E47').Value = '  +8.14%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue 'D48' '34.89'
$ws.Range('E48').Value = '  +2.71%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D49' '0.1207'
$ws.Range('E49').Value = '  +0.92%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D50' '0.05621'
$ws.Range('E50').Value = '  +2.21%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D51' '0.3855'
$ws.Range('E51').Value = '  +4.31%  '
